$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 1473.2858
$ws.Range("I40").Value = 1216
$ws.Range("J40").Value = 2116.5
$ws.Range("K40").Value = 1216
$ws.Range("L40").Value = 2116.5
$ws.Range("M40").Value = -1041
$ws.Range("N40").Value = -2466.5
$ws.Range("H41").Value = 302.94116
$ws.Range("J41").Value = 329.66666
$ws.Range("L41").Value = 329.66666
$ws.Range("N41").Value = -1209.66666
$ws.Range("H44").Value = 15025
$ws.Range("J44").Value = 15025
$ws.Range("L44").Value = 15025
$ws.Range("N44").Value = -15949
$ws.Range("H64").Value = 4405.6
$ws.Range("I64").Value = 4038.5
$ws.Range("K64").Value = 4038.5
$ws.Range("M64").Value = -3790.5
$ws.Range("H67").Value = 4405.6
$ws.Range("I67").Value = 4038.5
$ws.Range("K67").Value = 4038.5
$ws.Range("M67").Value = -3180.5
$ws.Range("H76").Value = 5558525.5
$ws.Range("I76").Value = 3200
$ws.Range("J76").Value = 13891514
$ws.Range("K76").Value = 3200
$ws.Range("L76").Value = 13891514
$ws.Range("M76").Value = -2885
$ws.Range("N76").Value = -13892144
$ws.Range("H79").Value = 5558525.5
$ws.Range("I79").Value = 3200
$ws.Range("J79").Value = 13891514
$ws.Range("K79").Value = 3200
$ws.Range("L79").Value = 13891514
$ws.Range("M79").Value = -2108
$ws.Range("N79").Value = -13893698
$ws.Range("H103").Value = 156638.75
$ws.Range("I103").Value = 208695.83
$ws.Range("J103").Value = 467.5
$ws.Range("K103").Value = 626087.49
$ws.Range("L103").Value = 1402.5
$ws.Range("M103").Value = -625501.49
$ws.Range("N103").Value = -2574.5
$ws.Range("H129").Value = 796.39
$ws.Range("I129").Value = 440.375
$ws.Range("J129").Value = 827.34784
$ws.Range("K129").Value = 1321.125
$ws.Range("L129").Value = 2482.04352
$ws.Range("M129").Value = 3678.875
$ws.Range("N129").Value = -12482.04352
$ws.Range("H132").Value = 2606.147
$ws.Range("I132").Value = 2687.6428
$ws.Range("K132").Value = 8062.928400000001
$ws.Range("M132").Value = -5532.928400000001

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H4").Value = 100
$ws.Range("I4").Value = 0
$ws.Range("K4").Value = 0
$ws.Range("M4").ClearContents()
$ws.Range("H39").Value = 4999
$ws.Range("I39").Value = 4999
$ws.Range("K39").Value = 4999
$ws.Range("M39").Value = -4479
$ws.Range("H42").Value = 16000
$ws.Range("J42").Value = 16000
$ws.Range("L42").Value = 16000
$ws.Range("N42").Value = -16972
$ws.Range("H44").Value = 23400
$ws.Range("J44").Value = 23400
$ws.Range("L44").Value = 23400
$ws.Range("N44").Value = -24376
$ws.Range("H63").Value = 1954936.5
$ws.Range("I63").Value = 1921.9231
$ws.Range("K63").Value = 1921.9231
$ws.Range("M63").Value = -1235.9231
$ws.Range("H66").Value = 1954936.5
$ws.Range("I66").Value = 1921.9231
$ws.Range("K66").Value = 9609.6155
$ws.Range("M66").Value = -6177.6155
$ws.Range("H132").Value = 18092.549
$ws.Range("I132").Value = 1648.8846
$ws.Range("K132").Value = 4946.6538
$ws.Range("M132").Value = -2416.6538

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H22").Value = 290.5
$ws.Range("I22").Value = 296.2857
$ws.Range("K22").Value = 296.2857
$ws.Range("M22").Value = -123.2857
$ws.Range("H86").Value = 1546.75
$ws.Range("I86").Value = 1375.25
$ws.Range("J86").Value = 2404.25
$ws.Range("K86").Value = 1375.25
$ws.Range("L86").Value = 2404.25
$ws.Range("M86").Value = -252.25
$ws.Range("N86").Value = -4650.25
$ws.Range("H89").Value = 1546.75
$ws.Range("I89").Value = 1375.25
$ws.Range("J89").Value = 2404.25
$ws.Range("K89").Value = 6876.25
$ws.Range("L89").Value = 12021.25
$ws.Range("M89").Value = -1260.25
$ws.Range("N89").Value = -23253.25
$ws.Range("H100").Value = 30817.2
$ws.Range("J100").Value = 30817.2
$ws.Range("L100").Value = 30817.2
$ws.Range("N100").Value = -32981.2
$ws.Range("H134").Value = 3681.4783
$ws.Range("I134").Value = 3889.238
$ws.Range("K134").Value = 11667.714
$ws.Range("M134").Value = -9132.714

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2943.718
$ws.Range("I31").Value = 2272.2856
$ws.Range("K31").Value = 2272.2856
$ws.Range("M31").Value = -1977.2856
$ws.Range("H34").Value = 2943.718
$ws.Range("I34").Value = 2272.2856
$ws.Range("K34").Value = 2272.2856
$ws.Range("M34").Value = -2070.2856
$ws.Range("H62").Value = 5836.1665
$ws.Range("I62").Value = 3502.5
$ws.Range("J62").Value = 7003
$ws.Range("K62").Value = 3502.5
$ws.Range("L62").Value = 7003
$ws.Range("M62").Value = -2878.5
$ws.Range("N62").Value = -8251
$ws.Range("H65").Value = 5836.1665
$ws.Range("I65").Value = 3502.5
$ws.Range("J65").Value = 7003
$ws.Range("K65").Value = 17512.5
$ws.Range("L65").Value = 35015
$ws.Range("M65").Value = -14392.5
$ws.Range("N65").Value = -41255
$ws.Range("H105").Value = 20833882
$ws.Range("I105").Value = 20833882
$ws.Range("K105").Value = 20833882
$ws.Range("M105").Value = -20832135
$ws.Range("H122").Value = 2322.5557
$ws.Range("I122").Value = 2322.5557
$ws.Range("K122").Value = 6967.6671
$ws.Range("M122").Value = -4517.6671
$ws.Range("H124").Value = 10111.111
$ws.Range("I124").Value = 9812.5
$ws.Range("J124").Value = 12500
$ws.Range("K124").Value = 9812.5
$ws.Range("L124").Value = 12500
$ws.Range("M124").Value = -7357.5
$ws.Range("N124").Value = -17410

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H36").Value = 2669.6
$ws.Range("I36").Value = 1914
$ws.Range("K36").Value = 5742
$ws.Range("M36").Value = -5573
$ws.Range("H68").Value = 100003
$ws.Range("I68").Value = 0
$ws.Range("J68").Value = 100003
$ws.Range("K68").Value = 0
$ws.Range("L68").Value = 300009
$ws.Range("M68").ClearContents()
$ws.Range("N68").Value = -301631
$ws.Range("H71").Value = 100003
$ws.Range("I71").Value = 0
$ws.Range("J71").Value = 100003
$ws.Range("K71").Value = 0
$ws.Range("L71").Value = 900027
$ws.Range("M71").ClearContents()
$ws.Range("N71").Value = -908139
$ws.Range("H131").Value = 777.03064
$ws.Range("J131").Value = 791.82794
$ws.Range("L131").Value = 2375.48382
$ws.Range("N131").Value = -12455.48382
$ws.Range("H141").Value = 3705.9285
$ws.Range("I141").Value = 3948.5
$ws.Range("J141").Value = 3382.5
$ws.Range("K141").Value = 11845.5
$ws.Range("L141").Value = 10147.5
$ws.Range("M141").Value = -6665.5
$ws.Range("N141").Value = -20507.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 2290.9167
$ws.Range("I97").Value = 754.125
$ws.Range("J97").Value = 5364.5
$ws.Range("K97").Value = 754.125
$ws.Range("L97").Value = 5364.5
$ws.Range("M97").Value = -258.125
$ws.Range("N97").Value = -6356.5
$ws.Range("H102").Value = 1779.5807
$ws.Range("I102").Value = 1780.2693
$ws.Range("K102").Value = 1780.2693
$ws.Range("M102").Value = -158.2692999999999
$ws.Range("H126").Value = 5264
$ws.Range("I126").Value = 4665.7393
$ws.Range("K126").Value = 13997.2179
$ws.Range("M126").Value = -11527.2179

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H127").Value = 35454.547
$ws.Range("J127").Value = 35454.547
$ws.Range("L127").Value = 35454.547
$ws.Range("N127").Value = -45374.547

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 3031722.8
$ws.Range("J107").Value = 5052216
$ws.Range("L107").Value = 15156648
$ws.Range("N107").Value = -15160488
